$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C1").Value = "deRegistered"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value = "NA"
}

$ws.Range("C1").Font.Bold = $true
